$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cash Deposit for 2021-07-16"
$ws.Range("B1").Value = "Credit Deposit for 2021-07-16"
$ws.Range("A2").Value = 520.52
$ws.Range("B2").Value = 437.06
